$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q1" sheet, positioned right before "总计"
#    (so the final tab order is 2021-Q2, 2021-Q3, 2022-Q1, 总计).
# ---------------------------------------------------------------------------
# NOTE: worksheet references in this host are positional (bound to an
# index), not stable object handles - they must be re-looked-up by name
# after any operation (Add/Move/Rename/Delete) that can shift indices.
$newSheet = $wb.Worksheets.Add()                 # lands at index 1
$totalSheet = $wb.Worksheets.Item("总计")         # look up AFTER the add
$newSheet.Move($totalSheet)                      # new sheet -> right before 总计

# Re-fetch the moved sheet by its now-settled position (3rd tab) and rename.
$q1Sheet = $wb.Worksheets.Item(3)
$q1Sheet.Name = "2022-Q1"

# Pull header/index formatting (bold font + border + centered, used by the
# sibling quarter sheets) from an existing sheet so the new one matches.
$fmtSrc = $wb.Worksheets.Item("2021-Q3")
$fmtSrc.Range("B1:H1").Copy($q1Sheet.Range("B1:H1"))
$fmtSrc.Range("A2").Copy($q1Sheet.Range("A2"))

# Header row text
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Data row (index column keeps the copied numeric/border style)
$q1Sheet.Range("A2").Value = 0

# The remaining fields must stay *text*, even though several of them look
# like plain numbers ("080005", "69.88", ...). Force text storage via a
# temporary "@" number format, then strip the format back off so no stray
# style index is left behind (matches the plain, un-styled cells we need).
$q1Sheet.Range("B2:G2").NumberFormat = "@"
$q1Sheet.Range("B2").Value = "080005"
$q1Sheet.Range("C2").Value = "长盛量化红利混合"
$q1Sheet.Range("D2").Value = "2.66"
$q1Sheet.Range("E2").Value = "69.88"
$q1Sheet.Range("F2").Value = "2.48"
$q1Sheet.Range("G2").Value = "0.0660"
$q1Sheet.Range("B2:G2").ClearFormats()

# Rank column is a genuine number
$q1Sheet.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: prepend a "2022-Q1" row, pushing the
#    previously-existing rows down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Grow the table by one row; copy the index-column (A) style down onto the
# new last row before it's populated, so every row keeps the same border.
$totalSheet.Range("A3").Copy($totalSheet.Range("A4"))

# Fill bottom-up with the final literal values so nothing is clobbered
# before it has been read.
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 8
$totalSheet.Range("D4").Value = 0.24

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q3"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.06

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.07000000000000001
